# preparation publication 0.2.0
# - bump Version 0.1.1 -> 0.2.0
# - bump Date 2023-10-20T07:19:33+00:00 -> 2023-10-20T08:59:58+00:00
# - insert a new "Jurisdiction" / "iso:code:3166:FR" row right after "Contact",
#   pushing Description/Purpose/Copyright/Immutable down by one row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3, column B)
$ws.Range("B3").Value = "0.2.0"

# Update Date value (row 8, column B)
$ws.Range("B8").Value = "2023-10-20T08:59:58+00:00"

# Row 15 is brand new (sheet used to end at row 14) - give it the same
# formatting as the rest of the data rows before we populate it
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# Shift rows 11-14 down to 12-15 (bottom-up so we don't clobber data we still need)
for ($r = 14; $r -ge 11; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value = $ws.Cells.Item($r, 2).Value2
}

# Insert the new Jurisdiction row at row 11
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
